$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 1.02),
    @("C2", 1.069504995555226),
    @("D2", 1.068227240103942),
    @("E2", 1.073200654573543),
    @("F2", 1.081748225534814),
    @("I2", 1.047873759189905),
    @("J2", 1.074439382078946),
    @("K2", 1.070932987304473),
    @("L2", 1.075893145079644),
    @("M2", 1.084418239958535),
    @("B3", 1.019999999999999),
    @("C3", 1.071079478604031),
    @("D3", 1.069436886922872),
    @("E3", 1.074584481136349),
    @("F3", 1.083168121359608),
    @("I3", 1.048271202457672),
    @("J3", 1.075668105253648),
    @("K3", 1.071957593096608),
    @("L3", 1.077092466562212),
    @("M3", 1.085655183955298),
    @("B4", 1.02),
    @("C4", 1.072096838609929),
    @("D4", 1.070218086496219),
    @("E4", 1.07547876324512),
    @("F4", 1.084085738166867),
    @("I4", 1.04852603424059),
    @("J4", 1.076461324540366),
    @("K4", 1.072618489562522),
    @("L4", 1.077866825694777),
    @("M4", 1.08645388699246),
    @("B5", 1.02),
    @("C5", 1.072524201142359),
    @("D5", 1.070546143554305),
    @("E5", 1.075854451335136),
    @("F5", 1.0844712351298),
    @("I5", 1.048632607461892),
    @("J5", 1.076794357895217),
    @("K5", 1.07289583460612),
    @("L5", 1.078191969401765),
    @("M5", 1.086789264817792),
    @("B6", 1.02),
    @("C6", 1.072595937746634),
    @("D6", 1.070601204865636),
    @("E6", 1.075917515478835),
    @("F6", 1.084535946167165),
    @("I6", 1.048650468915691),
    @("J6", 1.076850250267889),
    @("K6", 1.072942373128134),
    @("L6", 1.078246539356466),
    @("M6", 1.0868455531234),
    @("B7", 1.02),
    @("C7", 1.072102550357699),
    @("D7", 1.070222471415705),
    @("E7", 1.075483784253468),
    @("F7", 1.084090890245862),
    @("I7", 1.048527460467189),
    @("J7", 1.076465776254158),
    @("K7", 1.072622197402692),
    @("L7", 1.077871171833161),
    @("M7", 1.086458369879832),
    @("B8", 1.02),
    @("C8", 1.070037401668494),
    @("D8", 1.068636364429713),
    @("E8", 1.073668566370886),
    @("F8", 1.082228328075394),
    @("I8", 1.048008563015973),
    @("J8", 1.074855020889351),
    @("K8", 1.071279693808314),
    @("L8", 1.076298813183483),
    @("M8", 1.084836622864855),
    @("B9", 1.02),
    @("C9", 1.066386984248656),
    @("D9", 1.065829538792001),
    @("E9", 1.070460858742752),
    @("F9", 1.078937156792816),
    @("I9", 1.047076170300368),
    @("J9", 1.072002246235455),
    @("K9", 1.068897783719136),
    @("L9", 1.073514972110979),
    @("M9", 1.081965748771482),
    @("B10", 1.02),
    @("C10", 1.063945221332588),
    @("D10", 1.063949981410957),
    @("E10", 1.068315896080533),
    @("F10", 1.076736523396183),
    @("I10", 1.046442307633769),
    @("J10", 1.070090317454964),
    @("K10", 1.067298602830458),
    @("L10", 1.071649872603067),
    @("M10", 1.080042621920126),
    @("B11", 1.02),
    @("C11", 1.062885863023801),
    @("D11", 1.063134058149238),
    @("E11", 1.067385474704331),
    @("F11", 1.07578198843604),
    @("I11", 1.046164895434651),
    @("J11", 1.069259954167601),
    @("K11", 1.066603402380071),
    @("L11", 1.070839997973872),
    @("M11", 1.079207618948811),
    @("B12", 1.02),
    @("C12", 1.062492050905423),
    @("D12", 1.062830671409697),
    @("E12", 1.067039620850493),
    @("F12", 1.075427176476296),
    @("I12", 1.046061406748232),
    @("J12", 1.068951139315855),
    @("K12", 1.066344755399514),
    @("L12", 1.070538825601946),
    @("M12", 1.078897112415615),
    @("B13", 1.02),
    @("C13", 1.062576539550147),
    @("D13", 1.062895763325253),
    @("E13", 1.067113819292482),
    @("F13", 1.075503296520108),
    @("I13", 1.046083625627199),
    @("J13", 1.0690173985105),
    @("K13", 1.066400255097177),
    @("L13", 1.070603444001569),
    @("M13", 1.078963733042489),
    @("B14", 1.02),
    @("C14", 1.062853316962213),
    @("D14", 1.063108986602066),
    @("E14", 1.067356891542148),
    @("F14", 1.075752664819924),
    @("I14", 1.046156350131452),
    @("J14", 1.069234435250794),
    @("K14", 1.066582031115147),
    @("L14", 1.070815110130439),
    @("M14", 1.079181959548711),
    @("B15", 1.02),
    @("C15", 1.063023806054481),
    @("D15", 1.063240318359677),
    @("E15", 1.067506622495923),
    @("F15", 1.075906274810388),
    @("I15", 1.046201098985305),
    @("J15", 1.06936810804291),
    @("K15", 1.066693973665792),
    @("L15", 1.070945478172693),
    @("M15", 1.07931636960436),
    @("B16", 1.02),
    @("C16", 1.064015483228992),
    @("D16", 1.064004087530283),
    @("E16", 1.068377609901529),
    @("F16", 1.076799837385707),
    @("I16", 1.046460656261399),
    @("J16", 1.07014537292302),
    @("K16", 1.067344682638508),
    @("L16", 1.071703572786454),
    @("M16", 1.08009798973929),
    @("B17", 1.02),
    @("C17", 1.064636977903012),
    @("D17", 1.064482622888015),
    @("E17", 1.068923513174421),
    @("F17", 1.077359899327205),
    @("I17", 1.046622679042363),
    @("J17", 1.070632259905131),
    @("K17", 1.067752115997849),
    @("L17", 1.0721784915687),
    @("M17", 1.080587664851914),
    @("B18", 1.02),
    @("C18", 1.064999287258952),
    @("D18", 1.064761545915784),
    @("E18", 1.069241771766003),
    @("F18", 1.077686415684522),
    @("I18", 1.046716900192015),
    @("J18", 1.070916013541507),
    @("K18", 1.067989500596426),
    @("L18", 1.0724552847184),
    @("M18", 1.080873064919043),
    @("B19", 1.02),
    @("C19", 1.065122792106597),
    @("D19", 1.064856617982438),
    @("E19", 1.069350263231813),
    @("F19", 1.077797722710774),
    @("I19", 1.046748979085995),
    @("J19", 1.071012725835527),
    @("K19", 1.068070397953061),
    @("L19", 1.072549627109884),
    @("M19", 1.080970342014645),
    @("B20", 1.02),
    @("C20", 1.06457031795561),
    @("D20", 1.064431301158079),
    @("E20", 1.068864959272299),
    @("F20", 1.077299826411671),
    @("I20", 1.046605324927663),
    @("J20", 1.070580046346626),
    @("K20", 1.067708429646278),
    @("L20", 1.072127559976212),
    @("M20", 1.08053515006598),
    @("B21", 1.02),
    @("C21", 1.062771821806727),
    @("D21", 1.06304620645238),
    @("E21", 1.06728531989546),
    @("F21", 1.07567923917845),
    @("I21", 1.046134946897568),
    @("J21", 1.069170533917315),
    @("K21", 1.066528514213712),
    @("L21", 1.070752789419866),
    @("M21", 1.07911770698138),
    @("B22", 1.02),
    @("C22", 1.061639183684738),
    @("D22", 1.062173507868744),
    @("E22", 1.066290664878292),
    @("F22", 1.074658831363245),
    @("I22", 1.04583662289388),
    @("J22", 1.068282109645119),
    @("K22", 1.065784229649725),
    @("L22", 1.069886394517042),
    @("M22", 1.078224479959747),
    @("B23", 1.02),
    @("C23", 1.06223979549454),
    @("D23", 1.062636318027658),
    @("E23", 1.066818092562363),
    @("F23", 1.07519991180682),
    @("I23", 1.045995015501889),
    @("J23", 1.06875329197244),
    @("K23", 1.066179020928144),
    @("L23", 1.070345880935568),
    @("M23", 1.07869819082114),
    @("B24", 1.02),
    @("C24", 1.064600439328431),
    @("D24", 1.064454491852485),
    @("E24", 1.068891417747324),
    @("F24", 1.077326971265792),
    @("I24", 1.04661316738274),
    @("J24", 1.070603640145854),
    @("K24", 1.067728170445242),
    @("L24", 1.072150574447263),
    @("M24", 1.080558879914639),
    @("B25", 1.02),
    @("C25", 1.067332103586539),
    @("D25", 1.066556615598236),
    @("E25", 1.071291242741623),
    @("F25", 1.079789122874742),
    @("I25", 1.04731936762727),
    @("J25", 1.07274150493617),
    @("K25", 1.069515522435142),
    @("L25", 1.074236257133369),
    @("M25", 1.082709534673982)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}